# UC012 - Listar Liquidações Pendentes
# Change: swap the second-step text (Steps / Expected Results) between
# Test Case TC2 and Test Case TC3.
#
# Previously:
#   TC2 step 2 -> "Chefe Clica para realizar a liquidação."
#                 "SYSTEM Apresenta a tela de Registrar Liquidações"
#   TC3 step 2 -> "Chefe Clica para atribuir/desatribuir o registro a si mesmo."
#                 "SYSTEM Atualiza a lista de registros de solicitações, onde o nome
#                  deverá constar o nome do usuário logado (que se atribuiu como
#                  responsável pela liquidação) no campo de atribuição (no caso de
#                  desatribuição, o nome deverá ser removido)."
#
# After:
#   TC2 step 2 -> "Chefe Clica para atribuir/desatribuir o registro a si mesmo."
#                 "SYSTEM Atualiza a lista de registros de solicitações, onde o nome
#                  deverá constar o nome do usuário logado (que se atribuiu como
#                  responsável pela liquidação) no campo de atribuição (no caso de
#                  desatribuição, o nome deverá ser removido)."
#   TC3 step 2 -> "Chefe Clica para realizar a liquidação."
#                 "SYSTEM Apresenta a tela de Registrar Liquidações"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$attribText = "Chefe Clica para atribuir/desatribuir o registro a si mesmo."
$attribResult = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pela liquidação) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

$liquidText = "Chefe Clica para realizar a liquidação."
$liquidResult = "SYSTEM Apresenta a tela de Registrar Liquidações"

# TC2 is the block whose "Test Case ID" cell (column B) holds "TC2"; its
# second step row is two rows below the "#" header row, i.e. 5 rows below
# the Test Case ID row. In this workbook that is row 20.
$ws.Range("B20").Value = $attribText
$ws.Range("D20").Value = $attribResult

# TC3's second step row (row 28) gets the text that used to belong to TC2.
$ws.Range("B28").Value = $liquidText
$ws.Range("D28").Value = $liquidResult
